$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.035981583171684
$ws.Range("D2").Value = 1.043259295599308
$ws.Range("E2").Value = 1.034964889398478
$ws.Range("F2").Value = 1.050552412143832
$ws.Range("I2").Value = 1.04062089495411
$ws.Range("J2").Value = 1.041092390456456
$ws.Range("K2").Value = 1.046033665589697
$ws.Range("L2").Value = 1.037762851298454
$ws.Range("M2").Value = 1.053306371980609
$ws.Range("N2").Value = 1.005712725503983

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.036848356035097
$ws.Range("D3").Value = 1.043943365355206
$ws.Range("E3").Value = 1.035699518811472
$ws.Range("F3").Value = 1.051417456946365
$ws.Range("I3").Value = 1.04085231878188
$ws.Range("J3").Value = 1.041603354871948
$ws.Range("K3").Value = 1.046529066063774
$ws.Range("L3").Value = 1.038306976380271
$ws.Range("M3").Value = 1.053983750560193

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.037409729516643
$ws.Range("D4").Value = 1.044386415440684
$ws.Range("E4").Value = 1.036175693641938
$ws.Range("F4").Value = 1.051978085168057
$ws.Range("I4").Value = 1.041000964028836
$ws.Range("J4").Value = 1.041933844084615
$ws.Range("K4").Value = 1.046849355131163
$ws.Range("L4").Value = 1.038659214109823
$ws.Range("M4").Value = 1.054422303318094

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.037645852393098
$ws.Range("D5").Value = 1.044572770474701
$ws.Range("E5").Value = 1.0363760722721
$ws.Range("F5").Value = 1.052213983691403
$ws.Range("I5").Value = 1.041063190036541
$ws.Range("J5").Value = 1.042072747288374
$ws.Range("K5").Value = 1.04698393921221
$ws.Range("L5").Value = 1.038807330049266
$ws.Range("M5").Value = 1.054606727682175

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.037685505543328
$ws.Range("D6").Value = 1.04460406593154
$ws.Range("E6").Value = 1.036409728095723
$ws.Range("F6").Value = 1.052253604378661
$ws.Range("I6").Value = 1.041073622527081
$ws.Range("J6").Value = 1.042096067691758
$ws.Range("K6").Value = 1.047006532597863
$ws.Range("L6").Value = 1.03883220139372
$ws.Range("M6").Value = 1.054637696630249

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.037412884126626
$ws.Range("D7").Value = 1.044388905147592
$ws.Range("E7").Value = 1.036178370347966
$ws.Range("F7").Value = 1.051981236431107
$ws.Range("I7").Value = 1.041001796535876
$ws.Range("J7").Value = 1.041935700252935
$ws.Range("K7").Value = 1.046851153708945
$ws.Range("L7").Value = 1.038661193105571
$ws.Range("M7").Value = 1.054424767384247

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.036274405919201
$ws.Range("D8").Value = 1.043490394030603
$ws.Range("E8").Value = 1.035212989945837
$ws.Range("F8").Value = 1.050844573478875
$ws.Range("I8").Value = 1.040699333320647
$ws.Range("J8").Value = 1.041265101317679
$ws.Range("K8").Value = 1.04620114327733
$ws.Range("L8").Value = 1.037946708582462
$ws.Range("M8").Value = 1.053535243909375

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.034272259221553
$ws.Range("D9").Value = 1.041910325642142
$ws.Range("E9").Value = 1.033518224978971
$ws.Range("F9").Value = 1.048848495713123
$ws.Range("I9").Value = 1.04015795262465
$ws.Range("J9").Value = 1.040082409325876
$ws.Range("K9").Value = 1.045053742957521
$ws.Range("L9").Value = 1.036688919936836
$ws.Range("M9").Value = 1.051969721937102

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.032940269074439
$ws.Range("D10").Value = 1.040859212791866
$ws.Range("E10").Value = 1.032392757582808
$ws.Range("F10").Value = 1.047522500707305
$ws.Range("I10").Value = 1.039791428730079
$ws.Range("J10").Value = 1.039293339574343
$ws.Range("K10").Value = 1.044287539338563
$ws.Range("L10").Value = 1.035851294112001
$ws.Range("M10").Value = 1.050927429659646

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.032364178774991
$ws.Range("D11").Value = 1.040404627671755
$ws.Range("E11").Value = 1.031906476908553
$ws.Range("F11").Value = 1.046949472897795
$ws.Range("I11").Value = 1.039631400303726
$ws.Range("J11").Value = 1.038951534311405
$ws.Range("K11").Value = 1.043955478977021
$ws.Range("L11").Value = 1.035488822469965
$ws.Range("M11").Value = 1.050476453974149

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.032150295317042
$ws.Range("D12").Value = 1.040235859381374
$ws.Range("E12").Value = 1.031726010701302
$ws.Range("F12").Value = 1.046736797438589
$ws.Range("I12").Value = 1.039571760875346
$ws.Range("J12").Value = 1.038824553849998
$ws.Range("K12").Value = 1.043832094769758
$ws.Range("L12").Value = 1.035354219590471
$ws.Range("M12").Value = 1.0503089946838

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.03219616940014
$ws.Range("D13").Value = 1.040272056885354
$ws.Range("E13").Value = 1.031764714043582
$ws.Range("F13").Value = 1.04678240920795
$ws.Range("I13").Value = 1.039584562671625
$ws.Range("J13").Value = 1.038851792420161
$ws.Range("K13").Value = 1.043858562991032
$ws.Range("L13").Value = 1.035383090750168
$ws.Range("M13").Value = 1.050344912840161

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.032346497009902
$ws.Range("D14").Value = 1.040390675491326
$ws.Range("E14").Value = 1.031891556239931
$ws.Range("F14").Value = 1.046931889546963
$ws.Range("I14").Value = 1.039626474519612
$ws.Range("J14").Value = 1.038941038442818
$ws.Range("K14").Value = 1.043945280853761
$ws.Range("L14").Value = 1.03547769542785
$ws.Range("M14").Value = 1.050462610643784

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.032439132387907
$ws.Range("D15").Value = 1.040463771624831
$ws.Range("E15").Value = 1.031969729167744
$ws.Range("F15").Value = 1.047024012239966
$ws.Range("I15").Value = 1.039652271611499
$ws.Range("J15").Value = 1.038996023411598
$ws.Range("K15").Value = 1.043998705035564
$ws.Range("L15").Value = 1.035535989198768
$ws.Range("M15").Value = 1.050535135234496

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.032978516533659
$ws.Range("D16").Value = 1.04088939391795
$ws.Range("E16").Value = 1.032425052785472
$ws.Range("F16").Value = 1.047560554779647
$ws.Range("I16").Value = 1.039802021517811
$ws.Range("J16").Value = 1.039316021366594
$ws.Range("K16").Value = 1.044309571104812
$ws.Range("L16").Value = 1.035875355030293
$ws.Range("M16").Value = 1.050957366802307

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.03331703840869
$ws.Range("D17").Value = 1.041156524959712
$ws.Range("E17").Value = 1.0327109486945
$ws.Range("F17").Value = 1.047897419354643
$ws.Range("I17").Value = 1.039895602489991
$ws.Range("J17").Value = 1.039516712902177
$ws.Range("K17").Value = 1.044504492638904
$ws.Range("L17").Value = 1.036088291634855
$ws.Range("M17").Value = 1.051222314610032

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.033514557061471
$ws.Range("D18").Value = 1.04131239128647
$ws.Range("E18").Value = 1.032877808535117
$ws.Range("F18").Value = 1.048094016324164
$ws.Range("I18").Value = 1.039950059164907
$ws.Range("J18").Value = 1.039633759981075
$ws.Range("K18").Value = 1.044618159091393
$ws.Range("L18").Value = 1.036212515659936
$ws.Range("M18").Value = 1.051376887248903

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.033581916726351
$ws.Range("D19").Value = 1.041365546667061
$ws.Range("E19").Value = 1.032934720608604
$ws.Range("F19").Value = 1.048161069339287
$ws.Range("I19").Value = 1.039968605831062
$ws.Range("J19").Value = 1.039673667786167
$ws.Range("K19").Value = 1.044656911630709
$ws.Range("L19").Value = 1.036254876484653
$ws.Range("M19").Value = 1.051429598086597

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.033280711526039
$ws.Range("D20").Value = 1.041127858800844
$ws.Range("E20").Value = 1.032680264248442
$ws.Range("F20").Value = 1.047861265628302
$ws.Range("I20").Value = 1.039885575322235
$ws.Range("J20").Value = 1.039495181925264
$ws.Range("K20").Value = 1.044483582288246
$ws.Range("L20").Value = 1.036065443297878
$ws.Range("M20").Value = 1.051193884786678

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.032302226444815
$ws.Range("D21").Value = 1.040355742912228
$ws.Range("E21").Value = 1.031854199941622
$ws.Range("F21").Value = 1.046887866534261
$ws.Range("I21").Value = 1.039614137977768
$ws.Range("J21").Value = 1.038914758222056
$ws.Range("K21").Value = 1.043919745759685
$ws.Range("L21").Value = 1.035449835725399
$ws.Range("M21").Value = 1.050427950094546

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.031687605642652
$ws.Range("D22").Value = 1.039870774832716
$ws.Range("E22").Value = 1.031335747586075
$ws.Range("F22").Value = 1.046276851292288
$ws.Range("I22").Value = 1.039442330820629
$ws.Range("J22").Value = 1.038549714808013
$ws.Range("K22").Value = 1.043564995788028
$ws.Range("L22").Value = 1.035062983265977
$ws.Range("M22").Value = 1.049946684864849

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.032013371201165
$ws.Range("D23").Value = 1.040127818407932
$ws.Range("E23").Value = 1.031610500542395
$ws.Range("F23").Value = 1.046600666701552
$ws.Range("I23").Value = 1.039533517208406
$ws.Range("J23").Value = 1.038743241067248
$ws.Range("K23").Value = 1.043753078134824
$ws.Range("L23").Value = 1.035268041262519
$ws.Range("M23").Value = 1.050201782938827

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.033297125881278
$ws.Range("D24").Value = 1.041140811639603
$ws.Range("E24").Value = 1.032694128915124
$ws.Range("F24").Value = 1.047877601602702
$ws.Range("I24").Value = 1.039890106561835
$ws.Range("J24").Value = 1.039504910886932
$ws.Range("K24").Value = 1.044493030862231
$ws.Range("L24").Value = 1.036075767410948
$ws.Range("M24").Value = 1.051206730898115

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.03478937950228
$ws.Range("D25").Value = 1.042318419033518
$ws.Range("E25").Value = 1.033955598429503
$ws.Range("F25").Value = 1.049363704881868
$ws.Range("I25").Value = 1.040298903491952
$ws.Range("J25").Value = 1.04038827601514
$ws.Range("K25").Value = 1.045350602710319
$ws.Range("L25").Value = 1.037013935447002
$ws.Range("M25").Value = 1.05237420865781
